# Weekly data refresh: insert the newest week's Brócoli price record at the
# top of the data block (row 545), pushing all existing records down by one
# row. This mirrors the "Fruta / hortaliza, semanal" weekly ingestion job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 545; rows 545:634 shift down to 546:635 and the
# used range / dimension grows to A1:R635 automatically.
$ws.Rows("545:545").Insert()

$ws.Range("A545").Value = 7
$ws.Range("B545").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C545").Value = "Ñuble"
$ws.Range("D545").Value = 45218
$ws.Range("E545").Value = 16
$ws.Range("F545").Value = 100112023
$ws.Range("G545").Value = "Brócoli"
$ws.Range("H545").Value = "Sin especificar"
$ws.Range("I545").Value = "Primera"
$ws.Range("J545").Value = 300
$ws.Range("K545").Value = 1200
$ws.Range("L545").Value = 1200
$ws.Range("M545").Value = 1200
$ws.Range("N545").Value = "$/unidad"
$ws.Range("O545").Value = "Provincia de Diguillín"
$ws.Range("P545").Value = 1200
$ws.Range("Q545").Value = 1
$ws.Range("R545").Value = "Hortaliza"
